# The deck's "datetimeFigureOut" date field (the Header & Footer "Date and
# time" placeholder) was cached as 13.06.2024 on the Slide Master and on
# every one of its Custom Layouts. The presentation was re-saved on
# 31.07.2024, so every one of those cached date placeholders needs to show
# 31.07.2024 instead.
#
# ppPlaceholderDate = 16 reliably identifies the date placeholder shape
# (named "Datumsplatzhalter N" in this German template) on the master and
# on each layout, regardless of its shape index, so we use that instead of
# a hard-coded Shapes.Item(n).

$OldDate = "13.06.2024"
$NewDate = "31.07.2024"

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if (-not $isDatePlaceholder) {
            continue
        }

        if (-not $shp.HasTextFrame) {
            continue
        }

        $tr = $shp.TextFrame.TextRange
        $tr.Text = $NewDate
    }
}

# Slide Master.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom layout hanging off the master.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
